$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("D2").Value = 1185
$ws2.Range("U17").Value = 0.04
$ws2.Range("U18").Value = 0.02
$ws2.Range("U19").Value = 0.085
